$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order (player, position, team) - rows 2..18
$data = @(
    @("Gradey Dick",   "SG,SF",       "Toronto Raptors"),
    @("Josh Giddey",   "PG,SG,SF",    "Chicago Bulls"),
    @("Klay Thompson", "SG,SF",       "Dallas Mavericks"),
    @("Scottie Barnes","PG,SG,SF,PF", "Toronto Raptors"),
    @("Miles Bridges", "SF,PF",       "Charlotte Hornets"),
    @("Nikola Vucevic","PF,C",        "Chicago Bulls"),
    @("Mikal Bridges", "SG,SF,PF",    "New York Knicks"),
    @("Brook Lopez",   "C",           "Milwaukee Bucks"),
    @("Ja Morant",     "PG",          "Memphis Grizzlies"),
    @("Tyler Herro",   "PG,SG",       "Miami Heat"),
    @("Evan Mobley",   "PF,C",        "Cleveland Cavaliers"),
    @("De'Aaron Fox",  "PG",          "Sacramento Kings"),
    @("Shaedon Sharpe","SG,SF",       "Portland Trail Blazers"),
    @("DeMar DeRozan", "SF,PF",       "Sacramento Kings"),
    @("Kyle Kuzma",    "PF",          "Washington Wizards"),
    @("Luka Doncic",   "PG,SG",       "Dallas Mavericks"),
    @("Goga Bitadze",  "C",           "Orlando Magic")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
